# feat: Factory event registration with Reflection
#
# Adds a new "Is Default" boolean column (G) to the "Tile Types" sheet,
# right after the existing "Can Dispawn Humans" column (F).
# All tile types default to FALSE except "Grass" (row 8), which is TRUE.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tile Types")

# Carry over the same cell style used by the rest of the table (column F)
# onto the new column before writing values into it.
$ws.Range("F1:F8").Copy()
$ws.Range("G1:G8").PasteSpecial(-4122)

$ws.Range("G1").Value = "Is Default"

$ws.Range("G2").Value = $false
$ws.Range("G3").Value = $false
$ws.Range("G4").Value = $false
$ws.Range("G5").Value = $false
$ws.Range("G6").Value = $false
$ws.Range("G7").Value = $false
$ws.Range("G8").Value = $true
